$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E6: 650 -> 588
$ws.Range("E6").Value = 588

# Add new values for H6 and I6
$ws.Range("H6").Value = 650
$ws.Range("I6").Value = 12014
